$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bibi Cell Mundi
$ws.Range("D2").Value = 12213.72
$ws.Range("AG2").Value = 30754.27

# Row 3 - Bibi Cell Vieiralves
$ws.Range("D3").Value = 5466.8
$ws.Range("AG3").Value = 12385.05

# Row 4 - Bibi Cell Ponta Negra
$ws.Range("C4").Value = 2416
$ws.Range("D4").Value = 2403.81
$ws.Range("AG4").Value = 9354.82

# Row 5 - Bibi Cell Manauara
$ws.Range("C5").Value = 3763
$ws.Range("D5").Value = 2753
$ws.Range("AG5").Value = 9272

# Row 6 - total
$ws.Range("C6").Value = 21194.65
$ws.Range("D6").Value = 22837.33
$ws.Range("AG6").Value = 61766.14
